{"js": "// Replace the two-digit \u00f7 one-digit division problems/answers in the\n// worksheet table. Each non-empty table cell holds exactly one paragraph\n// / one run of the form \"NN\u00f7N=NN, N\"; we walk the cells in document\n// order (row-major) and swap each old expression for its new one,\n// matching the order the diff lists them in. We only touch the text of\n// the existing run (via the paragraph's Range + InsertLocation.replace)\n// so the run/paragraph formatting (font, size, justification, ...) is\n// preserved untouched - exactly like the diff, which only rewrites the\n// <w:t> content.\n\nconst oldToNew = [\n  [\"14\u00f75=2, 4\", \"93\u00f74=23, 1\"],\n  [\"48\u00f75=9, 3\", \"98\u00f77=14, 0\"],\n  [\"59\u00f72=29, 1\", \"97\u00f79=10, 7\"],\n  [\"65\u00f76=10, 5\", \"40\u00f79=4, 4\"],\n  [\"79\u00f75=15, 4\", \"68\u00f72=34, 0\"],\n  [\"47\u00f77=6, 5\", \"55\u00f74=13, 3\"],\n  [\"36\u00f79=4, 0\", \"85\u00f79=9, 4\"],\n  [\"16\u00f72=8, 0\", \"38\u00f76=6, 2\"],\n  [\"37\u00f77=5, 2\", \"20\u00f72=10, 0\"],\n  [\"74\u00f72=37, 0\", \"81\u00f74=20, 1\"],\n  [\"40\u00f74=10, 0\", \"49\u00f73=16, 1\"],\n  [\"31\u00f76=5, 1\", \"31\u00f75=6, 1\"],\n  [\"41\u00f75=8, 1\", \"28\u00f72=14, 0\"],\n  [\"23\u00f74=5, 3\", \"41\u00f73=13, 2\"],\n  [\"93\u00f77=13, 2\", \"32\u00f74=8, 0\"],\n  [\"96\u00f79=10, 6\", \"69\u00f72=34, 1\"],\n  [\"51\u00f73=17, 0\", \"84\u00f73=28, 0\"],\n  [\"46\u00f74=11, 2\", \"51\u00f75=10, 1\"],\n  [\"22\u00f74=5, 2\", \"94\u00f72=47, 0\"],\n  [\"23\u00f74=5, 3\", \"93\u00f79=10, 3\"],\n  [\"56\u00f79=6, 2\", \"73\u00f72=36, 1\"],\n  [\"17\u00f78=2, 1\", \"67\u00f74=16, 3\"],\n  [\"16\u00f74=4, 0\", \"24\u00f76=4, 0\"],\n  [\"39\u00f75=7, 4\", \"28\u00f76=4, 4\"],\n  [\"90\u00f78=11, 2\", \"42\u00f73=14, 0\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst colCount = table.values.length > 0 ? table.values[0].length : 0;\n\n// Collect every non-empty cell (row-major order) so we can line it up\n// with the diff's sequence of old -> new replacements.\nconst cellCoords = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const text = (table.values[r][c] || \"\").trim();\n    if (text.length > 0) {\n      cellCoords.push([r, c]);\n    }\n  }\n}\n\nlet mapIndex = 0;\nfor (let i = 0; i < cellCoords.length && mapIndex < oldToNew.length; i++) {\n  const [r, c] = cellCoords[i];\n  const text = table.values[r][c].trim();\n  const [oldText, newText] = oldToNew[mapIndex];\n  if (text !== oldText) {\n    throw new Error(\n      `Cell (${r},${c}) text \"${text}\" did not match expected \"${oldText}\" (map index ${mapIndex}).`\n    );\n  }\n  const cell = table.getCell(r, c);\n  const para = cell.body.paragraphs.getFirst();\n  const rng = para.getRange();\n  rng.insertText(newText, Word.InsertLocation.replace);\n  mapIndex++;\n}\n\nawait context.sync();\n\nif (mapIndex !== oldToNew.length) {\n  throw new Error(\n    `Only replaced ${mapIndex} of ${oldToNew.length} expected cells.`\n  );\n}\n", "ps1": "# Replace the two-digit / one-digit division problems/answers in the\n# worksheet table. Each non-empty table cell holds exactly one paragraph\n# / one run of the form \"NN\u00f7N=NN, N\"; we walk the cells in document\n# order (row-major, 1-based COM indices) and swap each old expression\n# for its new one, matching the order the diff lists them in.\n#\n# We assign text directly to the existing Cell.Range (not the whole\n# cell/paragraph object), so the run/paragraph formatting (font, size,\n# justification, ...) is preserved untouched - exactly like the diff,\n# which only rewrites the <w:t> content.\n\n$oldToNew = @(\n    @(\"14\u00f75=2, 4\", \"93\u00f74=23, 1\"),\n    @(\"48\u00f75=9, 3\", \"98\u00f77=14, 0\"),\n    @(\"59\u00f72=29, 1\", \"97\u00f79=10, 7\"),\n    @(\"65\u00f76=10, 5\", \"40\u00f79=4, 4\"),\n    @(\"79\u00f75=15, 4\", \"68\u00f72=34, 0\"),\n    @(\"47\u00f77=6, 5\", \"55\u00f74=13, 3\"),\n    @(\"36\u00f79=4, 0\", \"85\u00f79=9, 4\"),\n    @(\"16\u00f72=8, 0\", \"38\u00f76=6, 2\"),\n    @(\"37\u00f77=5, 2\", \"20\u00f72=10, 0\"),\n    @(\"74\u00f72=37, 0\", \"81\u00f74=20, 1\"),\n    @(\"40\u00f74=10, 0\", \"49\u00f73=16, 1\"),\n    @(\"31\u00f76=5, 1\", \"31\u00f75=6, 1\"),\n    @(\"41\u00f75=8, 1\", \"28\u00f72=14, 0\"),\n    @(\"23\u00f74=5, 3\", \"41\u00f73=13, 2\"),\n    @(\"93\u00f77=13, 2\", \"32\u00f74=8, 0\"),\n    @(\"96\u00f79=10, 6\", \"69\u00f72=34, 1\"),\n    @(\"51\u00f73=17, 0\", \"84\u00f73=28, 0\"),\n    @(\"46\u00f74=11, 2\", \"51\u00f75=10, 1\"),\n    @(\"22\u00f74=5, 2\", \"94\u00f72=47, 0\"),\n    @(\"23\u00f74=5, 3\", \"93\u00f79=10, 3\"),\n    @(\"56\u00f79=6, 2\", \"73\u00f72=36, 1\"),\n    @(\"17\u00f78=2, 1\", \"67\u00f74=16, 3\"),\n    @(\"16\u00f74=4, 0\", \"24\u00f76=4, 0\"),\n    @(\"39\u00f75=7, 4\", \"28\u00f76=4, 4\"),\n    @(\"90\u00f78=11, 2\", \"42\u00f73=14, 0\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$mapIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($mapIndex -ge $oldToNew.Length) {\n            break\n        }\n        $cell = $t.Cell($r, $c)\n        $text = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($text.Length -eq 0) {\n            continue\n        }\n        $pair = $oldToNew[$mapIndex]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        if ($text -ne $oldText) {\n            throw \"Cell ($r,$c) text '$text' did not match expected '$oldText' (map index $mapIndex).\"\n        }\n        $cell.Range.Text = $newText\n        $mapIndex++\n    }\n}\n\nif ($mapIndex -ne $oldToNew.Length) {\n    throw \"Only replaced $mapIndex of $($oldToNew.Length) expected cells.\"\n}\n\nWrite-Output \"Replaced $mapIndex cells.\"\n"}
